$d = $word.ActiveDocument

# Map of old text -> new text, taken from the diff.
$replacements = [ordered]@{
    "2024-05-20 Monday" = "2024-05-21 Tuesday"
    "430÷8=53, 6"        = "386÷4=96, 2"
    "152÷3=50, 2"        = "904÷7=129, 1"
    "314÷6=52, 2"        = "860÷4=215, 0"
    "715÷2=357, 1"       = "894÷5=178, 4"
    "442÷6=73, 4"        = "503÷3=167, 2"
    "378÷7=54, 0"        = "395÷8=49, 3"
    "931÷9=103, 4"       = "229÷8=28, 5"
    "981÷6=163, 3"       = "621÷9=69, 0"
    "684÷8=85, 4"        = "235÷7=33, 4"
    "464÷7=66, 2"        = "258÷5=51, 3"
    "573÷7=81, 6"        = "399÷3=133, 0"
    "554÷4=138, 2"       = "652÷5=130, 2"
    "497÷9=55, 2"        = "722÷5=144, 2"
    "930÷2=465, 0"       = "658÷3=219, 1"
    "643÷3=214, 1"       = "747÷2=373, 1"
    "745÷6=124, 1"       = "212÷6=35, 2"
    "495÷3=165, 0"       = "136÷8=17, 0"
    "404÷7=57, 5"        = "908÷5=181, 3"
    "900÷6=150, 0"       = "546÷8=68, 2"
    "267÷6=44, 3"        = "437÷5=87, 2"
    "581÷9=64, 5"        = "340÷6=56, 4"
    "489÷9=54, 3"        = "578÷7=82, 4"
    "691÷6=115, 1"       = "443÷3=147, 2"
    "547÷6=91, 1"        = "502÷9=55, 7"
    "702÷3=234, 0"       = "657÷8=82, 1"
}

foreach ($old in $replacements.Keys) {
    $new = $replacements[$old]
    $range = $d.Content
    [void]$range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "Replacements applied."
